$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows continuing the series through 2021-09-01 (dates as Excel
# serials), copying column A's existing date style (row 357) forward so the
# new cells keep the same number format / border / font instead of picking
# up Excel's bare default style.
$data = @(
    @(358, 44432, 0, 8, 269.7235333782872),
    @(359, 44433, 0, 6, 202.2926500337155),
    @(360, 44434, 0, 3, 101.1463250168577),
    @(361, 44435, 0, 2, 67.43088334457181),
    @(362, 44436, 1, 3, 101.1463250168577),
    @(363, 44437, 1, 4, 134.8617666891436),
    @(364, 44438, 0, 2, 67.43088334457181),
    @(365, 44439, 3, 5, 168.5772083614295),
    @(366, 44440, 0, 5, 168.5772083614295)
)

foreach ($entry in $data) {
    $r = $entry[0]

    # Copy A357's style into the new A cell (carries the date number
    # format / border / font), then overwrite with the real date value.
    $ws.Cells.Item(357, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $entry[1]

    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
